$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J: header "target" (matching the bold/centered/wrap style used
# by the other header cells A1:H1), then "h" for every data row.
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "target"

$ws.Range("J2:J6").Value = "h"

# Move / reflect the active selection as it was left after the edit.
[void]$ws.Range("J9").Select()
